# Updated code with user registration
# Replace the sample transaction values in row 7 with a new registered
# transaction's data (new shared strings get appended automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "02-19-2019 18:57:57"
$ws.Range("B7").Value = "FT19021900040"

# C7 holds a long numeric-looking transaction id that must stay text
# (matches the original shared-string cell type), so prefix it with an
# apostrophe -- Excel's standard "force text" entry syntax -- to avoid
# it being coerced into a floating point number. Restore the cell style
# afterwards so the quote-prefix formatting doesn't change the cell's
# look (matches original, which carried no explicit style).
$ws.Range("C7").Value = "'20190219041910569"
$ws.Range("C7").Style = "Normal"

$ws.Range("G7").Value = "Automation FIftyfour"
